# Post-import cleanup: drop the temporary staging sheet ("Worksheet 1")
# that was used while importing data into tables b01/b02/b04, clear the
# leftover header-index row (row 4) on the report sheet, and leave the
# cursor on A5 where the real data now starts.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$ws = $wb.Worksheets.Item("Worksheet")

# Row 4 held a throwaway 0..10 column-index row used during the import;
# clear it so it disappears from the saved sheet.
$ws.Range("A4:K4").ClearContents()

# Remove the scratch "Worksheet 1" tab left over from the import process.
$wb.Worksheets.Item("Worksheet 1").Delete()

# Make sure the remaining sheet is active and leave the selection on A5,
# the first cell of the real data that now follows immediately after
# row 3.
$ws.Activate()
$ws.Range("A5").Select()
